$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.199797511100769
$ws.Range("B1").Value = 2.398146867752075
$ws.Range("C1").Value = 4.475554943084717
$ws.Range("D1").Value = 2.653626918792725
$ws.Range("E1").Value = 1.105576515197754
